# ---------------------------------------------------------------------------
# Fix: Corrigir cálculo de credenciamentos na exportação Excel
#
# Updates the weekly report workbook:
#   - GN "Cristian" (Alegrete) replaced by "Renan" (AG 003 - Zona Norte)
#   - Refreshed totals / percentages on all 4 sheets
#   - Per-GN detail sheet (tab 2) detail tables rebuilt for the new GN
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-Num($ws, $row, $col, $val) {
    # Plain numeric value -> stored as a Number cell (matches original type).
    $ws.Cells.Item($row, $col).Value = $val
}

function Set-Txt($ws, $row, $col, $val) {
    # Force text storage even when $val looks numeric/percent/date-ish, by
    # using the classic leading-apostrophe "text prefix" - otherwise Excel's
    # COM layer auto-coerces strings like "160.0%" or "4.00" into numbers.
    $ws.Cells.Item($row, $col).Value = "'" + $val
}

# ---------------------------------------------------------------------------
# 1) Rename the GN tab: Cristian -> Renan
# ---------------------------------------------------------------------------
$wsGn = $wb.Worksheets.Item(2)
$wsGn.Name = "Renan"

# ---------------------------------------------------------------------------
# 2) Sheet 1: "Resumo Executivo"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

Set-Txt $ws1 5 2 "13:03:39"

Set-Num $ws1 8 2 30
Set-Txt $ws1 9 2 "R$ 73.491,33"
Set-Num $ws1 10 2 52
Set-Num $ws1 11 2 69
Set-Num $ws1 12 2 16

Set-Txt $ws1 18 1 "Renan"
Set-Txt $ws1 18 2 "AG 003 - Zona Norte"
Set-Num $ws1 18 3 16
Set-Txt $ws1 18 5 "160.0%"
Set-Txt $ws1 18 6 "R$ 33.506,11"
Set-Txt $ws1 18 8 "22.3%"
Set-Num $ws1 18 9 19
Set-Txt $ws1 18 11 "63.3%"
Set-Num $ws1 18 13 4
Set-Txt $ws1 18 14 "57.1%"

# ---------------------------------------------------------------------------
# 3) Sheet 2: "Renan" (was "Cristian") - per-GN detail
# ---------------------------------------------------------------------------
Set-Txt $wsGn 1 1 "📋 DETALHAMENTO COMPLETO - RENAN"
Set-Txt $wsGn 2 2 "AG 003 - Zona Norte"

Set-Num $wsGn 7 2 16
Set-Txt $wsGn 7 4 "160.0%"
Set-Txt $wsGn 7 5 "META BATIDA"

Set-Txt $wsGn 8 2 "R$ 33.506,11"
Set-Txt $wsGn 8 4 "22.3%"

Set-Num $wsGn 9 2 19
Set-Txt $wsGn 9 4 "63.3%"

Set-Num $wsGn 12 2 4
Set-Num $wsGn 13 2 7
Set-Txt $wsGn 14 2 "57.1%"
Set-Txt $wsGn 15 2 "4.00"
Set-Txt $wsGn 16 2 "4.75"
Set-Num $wsGn 17 2 22
Set-Num $wsGn 18 2 0

# --- Rebuild the tail of the sheet (rows 24-40 -> rows 24-37) -------------
# The per-day detail / accreditation / simulation / PJ-manager rows for the
# previous GN no longer apply, and one fewer "days with accreditations" row
# is present, so the whole block collapses from 17 rows down to 14.
$wsGn.Range("A24:J40").ClearContents()

Set-Txt $wsGn 24 1 ""

Set-Txt $wsGn 25 1 "📝 DETALHAMENTO COMPLETO DOS CREDENCIAMENTOS"

Set-Txt $wsGn 26 1 "Data"
Set-Txt $wsGn 26 2 "Dia da Semana"
Set-Txt $wsGn 26 3 "EC"
Set-Txt $wsGn 26 4 "Volume (R$)"
Set-Txt $wsGn 26 5 "RA"
Set-Txt $wsGn 26 6 "Qual Oferta?"
Set-Txt $wsGn 26 7 "Instala Direto"
Set-Txt $wsGn 26 8 "Gerente PJ"
Set-Txt $wsGn 26 9 "Horário"

Set-Txt $wsGn 27 1 ""

Set-Txt $wsGn 28 1 "🔍 DETALHAMENTO COMPLETO DAS SIMULAÇÕES"

Set-Txt $wsGn 29 1 "Data"
Set-Txt $wsGn 29 2 "Dia da Semana"
Set-Txt $wsGn 29 3 "CNPJ"
Set-Txt $wsGn 29 4 "Empresa"
Set-Txt $wsGn 29 5 "Faturamento (R$)"
Set-Txt $wsGn 29 6 "Comentários"
Set-Txt $wsGn 29 7 "Horário"

Set-Txt $wsGn 30 1 ""

Set-Txt $wsGn 31 1 "👥 ANÁLISE DOS GERENTES PJ"

Set-Txt $wsGn 32 1 "Gerente PJ"
Set-Txt $wsGn 32 2 "Total Credenciamentos"
Set-Txt $wsGn 32 3 "Total Volume (R$)"

Set-Txt $wsGn 33 1 ""

Set-Txt $wsGn 34 1 "📊 RESUMO SEMANAL AVANÇADO"

Set-Txt $wsGn 35 1 "Dias com Credenciamentos:"
Set-Num $wsGn 35 2 4

Set-Txt $wsGn 36 1 "Dias com Simulações:"
Set-Num $wsGn 36 2 0

Set-Txt $wsGn 37 1 "Gerentes PJ Envolvidos:"
Set-Txt $wsGn 37 2 ""

# ---------------------------------------------------------------------------
# 4) Sheet 3: "Análise Comparativa"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

Set-Txt $ws3 6 2 "Renan"
Set-Txt $ws3 6 3 "AG 003 - Zona Norte"
Set-Txt $ws3 6 4 "16"
Set-Txt $ws3 6 6 "160.0%"
Set-Txt $ws3 6 7 "META BATIDA"
Set-Txt $ws3 6 8 "4"

Set-Txt $ws3 10 2 "Renan"
Set-Txt $ws3 10 3 "AG 003 - Zona Norte"
Set-Txt $ws3 10 4 "R$ 33.506,11"
Set-Txt $ws3 10 6 "22.3%"
Set-Txt $ws3 10 8 "4"

Set-Txt $ws3 14 1 "Renan"
Set-Txt $ws3 14 2 "AG 003 - Zona Norte"
Set-Txt $ws3 14 3 "4"
Set-Txt $ws3 14 4 "57.1%"
Set-Txt $ws3 14 5 "4.00"
Set-Txt $ws3 14 6 "4.75"

# ---------------------------------------------------------------------------
# 5) Sheet 4: "Cronograma Semanal"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

Set-Txt $ws4 6 1 "Renan"
Set-Txt $ws4 6 2 "-"
